$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "26.439.70"; E = "  +0.14%  " }
    @{ Row = 3; D = "1.810.64"; E = "  +0.54%  " }
    @{ Row = 4; D = "1.007"; E = "  -0.05%  " }
    @{ Row = 5; D = "1.006"; E = $null }
    @{ Row = 6; D = "306.13"; E = $null }
    @{ Row = 7; D = "0.4499"; E = "  -0.63%  " }
    @{ Row = 8; D = "0.3581"; E = "  -1.62%  " }
    @{ Row = 9; D = "46.39"; E = "  +3.31%  " }
    @{ Row = 10; D = "0.07054"; E = "  -0.38%  " }
    @{ Row = 11; D = "0.8890"; E = "  +2.23%  " }
    @{ Row = 12; D = "0.07792"; E = "  +0.53%  " }
    @{ Row = 13; D = "19.33"; E = "  +0.48%  " }
    @{ Row = 14; D = "1.788.00"; E = "  -0.70%  " }
    @{ Row = 15; D = "5.269"; E = "  +0.56%  " }
    @{ Row = 16; D = "6.303"; E = "  -0.22%  " }
    @{ Row = 17; D = "84.57"; E = "  -1.44%  " }
    @{ Row = 18; D = "1.009"; E = "  +0.04%  " }
    @{ Row = 19; D = "0.000008523"; E = "  -0.18%  " }
    @{ Row = 20; D = "1.007"; E = "  -0.01%  " }
    @{ Row = 21; D = "26.471.20"; E = "  +0.11%  " }
    @{ Row = 22; D = "14.16"; E = "  -0.09%  " }
    @{ Row = 23; D = "4.961"; E = "  +0.25%  " }
    @{ Row = 24; D = "2.034.10"; E = "  +0.60%  " }
    @{ Row = 25; D = "10.51"; E = "  +1.47%  " }
    @{ Row = 26; D = "1.949"; E = "  -0.66%  " }
    @{ Row = 27; D = "151.62"; E = "  +0.89%  " }
    @{ Row = 28; D = "17.77"; E = "  -0.32%  " }
    @{ Row = 29; D = "2.055"; E = "  +4.05%  " }
    @{ Row = 30; D = "112.10"; E = "  -0.46%  " }
    @{ Row = 31; D = "4.854"; E = "  +0.44%  " }
    @{ Row = 32; D = "0.08691"; E = "  +0.61%  " }
    @{ Row = 33; D = "3.115"; E = "  +2.91%  " }
    @{ Row = 34; D = "2.770"; E = "  +10.76%  " }
    @{ Row = 35; D = "0.7361"; E = "  +1.43%  " }
    @{ Row = 36; D = "4.443"; E = "  +0.39%  " }
    @{ Row = 37; D = "1.109"; E = "  +0.44%  " }
    @{ Row = 38; D = "1.073"; E = "  +0.21%  " }
    @{ Row = 39; D = "0.01925"; E = "  +0.35%  " }
    @{ Row = 40; D = "0.05126"; E = "  +1.68%  " }
    @{ Row = 41; D = "2.898"; E = "  +1.04%  " }
    @{ Row = 42; D = "0.5089"; E = "  +3.66%  " }
    @{ Row = 43; D = "6.767"; E = "  -2.17%  " }
    @{ Row = 44; D = "0.1505"; E = "  -3.65%  " }
    @{ Row = 45; D = "8.043"; E = "  -0.39%  " }
    @{ Row = 46; D = "0.4669"; E = "  +2.00%  " }
    @{ Row = 47; D = "1.007"; E = "  +0.11%  " }
    @{ Row = 48; D = "9.931"; E = "  +0.52%  " }
    @{ Row = 49; D = "99.85"; E = "  -1.29%  " }
    @{ Row = 50; D = "1.569"; E = "  -0.31%  " }
    @{ Row = 51; D = "0.05999"; E = "  +0.36%  " }
)

foreach ($u in $updates) {
    $dCell = $ws.Cells.Item($u.Row, 4)
    $origStyle = $dCell.Style
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.Style = $origStyle
    if ($u.E -ne $null) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

$wb.Save()